$p = $ppt.ActivePresentation

# --- Slide 1: Title - "Case Study 1: Beers and Brews" -> split into two runs,
#     second run becomes "BrewERIEs" ---
$s1 = $p.Slides.Item(1)
$titleShape = $s1.Shapes.Item(1)
$titleRange = $titleShape.TextFrame.TextRange
$titleText = $titleRange.Text
$brewsIndex = $titleText.IndexOf("Brews") + 1
$brewsRange = $titleRange.Characters($brewsIndex, 5)
$brewsRange.Text = "BrewERIEs"

# --- Slide 3: "AV data" -> "ABV data" ---
$s3 = $p.Slides.Item(3)
$s3.Shapes.Item(3).TextFrame.TextRange.Text = "Almost half (~42%) of IBU data was missing, and a little over 2.5% of ABV data was missing."

# --- Slide 4: "applying means of IBU and ABV" -> "applying average IBU and ABV" ---
$s4 = $p.Slides.Item(4)
$s4.Shapes.Item(3).TextFrame.TextRange.Text = "Filled missing values by applying average IBU and ABV per style of beer for greater accuracy."

# --- Slide 7: second run text update about the right skew ---
$s7 = $p.Slides.Item(7)
$s7TextRange = $s7.Shapes.Item(3).TextFrame.TextRange
$s7TextRange.Runs(2).Text = "One thing to note in particular, is the right skew in the distribution with the mean hovering around 5% ABV. This is a reflection of market competition, state laws, and history regarding German purity laws."

# --- Slide 8: delete the leftover "Rectangle 1" shape ---
$s8 = $p.Slides.Item(8)
$s8.Shapes.Item("Rectangle 1").Delete()
